# Weekly update: insert a new price record as row 174, shifting the
# existing rows 174:281 down to 175:282 (dimension grows to A1:T282).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(174).Insert()

$ws.Range("A174").Value = 7
$ws.Range("B174").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C174").Value = "Ñuble"
$ws.Range("D174").Value = 44488
$ws.Range("E174").Value = 16
$ws.Range("F174").Value = "Fruta"
$ws.Range("G174").Value = 100102
$ws.Range("H174").Value = "Cítricos"
$ws.Range("I174").Value = 100102005
$ws.Range("J174").Value = "Naranja"
$ws.Range("K174").Value = "Navel Late"
$ws.Range("L174").Value = "Primera"
$ws.Range("M174").Value = 300
$ws.Range("N174").Value = 7500
$ws.Range("O174").Value = 8000
$ws.Range("P174").Value = 7750
$ws.Range("Q174").Value = "$/bandeja 15 kilos granel"
$ws.Range("R174").Value = "Región de O'Higgins"
$ws.Range("S174").Value = 517
$ws.Range("T174").Value = 15
